$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CNPJ values in B2 and B3
$ws.Range("B2").Value = 11008634000107
$ws.Range("B3").Value = 11008634000107

# Update the active cell / selection on the sheet
$ws.Range("A12").Select() | Out-Null
